$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing cell B21 (week 20) with the new case count
$ws.Range("B21").Value = 305

# Add new row for week 21
$ws.Range("A22").Value = 21
$ws.Range("B22").Value = 6
